# Fixed tests for Error Messages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant for top vertical alignment (xlVAlignTop)
$xlVAlignTop = -4160

# Green fill color used elsewhere in the sheet (matches theme fill used for D/E columns)
$greenColor = $ws.Range("D4").Interior.Color

# --- Column widths (F and G split into two differently-sized columns) ---
$ws.Columns.Item(6).ColumnWidth = 18
$ws.Columns.Item(7).ColumnWidth = 26.453125

# --- Header row 3: "Error message" -> "Error messages" ---
$ws.Range("G3").Value = "Error messages"

# --- Row 4 ---
$ws.Rows.Item(4).RowHeight = 29
$ws.Range("A4:C4").VerticalAlignment = $xlVAlignTop
$ws.Range("D4:E4").VerticalAlignment = $xlVAlignTop
$ws.Range("F4").Value = "Check Page 2: some values, full JSON"
$ws.Range("F4").Interior.Color = $greenColor
$ws.Range("F4").WrapText = $true

# --- Row 5 ---
$ws.Range("F5").Value = "Check Full JSON"
$ws.Range("F5").Interior.Color = $greenColor

# --- Row 6 ---
$ws.Range("G6").ClearContents() | Out-Null

# --- Row 8 ---
$ws.Range("F8").Value = "Check Full JSON"
$ws.Range("F8").Interior.Color = $greenColor

# --- Row 9 ---
$ws.Range("G9").ClearContents() | Out-Null

# --- Row 15 ---
$ws.Rows.Item(15).RowHeight = 29
$ws.Range("A15:C15").VerticalAlignment = $xlVAlignTop
$ws.Range("D15").VerticalAlignment = $xlVAlignTop
$ws.Range("E15:F15").VerticalAlignment = $xlVAlignTop
$ws.Range("G15").Value = "Errors: User without email, user without password"
$ws.Range("G15").Interior.Color = $greenColor
$ws.Range("G15").VerticalAlignment = $xlVAlignTop
$ws.Range("G15").WrapText = $true

# --- Row 17 ---
$ws.Range("G17").Value = "Error: User without email"
$ws.Range("G17").Interior.Color = $greenColor

# --- Selection / view ---
$ws.Range("F13").Select() | Out-Null
